# Fill in the "Acceptance" column (D) for all review rows with "Accepted",
# matching the commit that records these review points as accepted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D13").Value = "Accepted"

# Update the saved selection to reflect where the user left the cursor.
[void]$ws.Range("D7").Select()
